# Apply updated quarterly values to the MCHP balance sheet worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MCHP")

# Row 4: Inventory
$ws.Range("B4").Value = 666000000.0
$ws.Range("C4").Value = 661000000.0
$ws.Range("D4").Value = 657000000.0
$ws.Range("E4").Value = 686000000.0
$ws.Range("F4").Value = 709000000.0

# Row 15: Accounts Payable
$ws.Range("B15").Value = 260000000.0
$ws.Range("C15").Value = 217000000.0
$ws.Range("D15").Value = 211000000.0
$ws.Range("E15").Value = 247000000.0
$ws.Range("F15").Value = 209000000.0

# Row 22: Long Term Tax Liability (Deferred)
$ws.Range("B22").Value = -1619000000.0
$ws.Range("C22").Value = -1553000000.0
$ws.Range("D22").Value = -1527000000.0
$ws.Range("E22").Value = -1430000000.0
$ws.Range("F22").Value = -1352000000.0
